$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "331.77"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.84%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.70"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.42%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.548"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.49%"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.68%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.056"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.27%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9770"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.52%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.1124"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-3.36%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1905"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.77%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "10.23"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-13.84%"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.46%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04698"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.16%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1059"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.90%"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-2.11%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.04107"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-3.08%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005930"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.79%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.349"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.66%"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.20%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.20%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-3.53%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1383"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.89%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2570"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.48%"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "3.72%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004396"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.62%"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "7.26%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003735"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-6.12%"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02770"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "5.14%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05735"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "3.40%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007615"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.56%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1423"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.15%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007527"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-7.13%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001955"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008274"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-6.69%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00007031"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-4.05%"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.27%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0005793"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.32%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003567"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-0.26%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002517"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "9.35%"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.27%"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.27%"
